# Weekly "symbol list" refresh of the cryptos sheet (GitHub Actions scrape).
# GateToken jumped in rank from row 17 up to row 6, pushing FTXToken..LEO
# down by one row each (their Coin/Link/Price/Volume move along with them),
# and every other row's Price (D) / Volume(1h) (E) got refreshed in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, then the new values for columns B (Coin), C (Link),
# D (Price) and E (Volume(1h)) that actually changed for that row. Columns
# left as $null are untouched for that row.
$updates = @(
    @{ Row = 2;  D = '305.91';  E = '2.79%' }
    @{ Row = 3;  D = '44.17' }
    @{ Row = 4;  D = '5.098';   E = '1.44%' }
    @{ Row = 5;  D = '0.07929'; E = '5.04%' }
    @{ Row = 6;  B = 'GateToken';                         C = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt';                         D = '4.443';      E = '1.61%' }
    @{ Row = 7;  B = 'FTXToken';                          C = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt';                            D = '1.595';      E = '0.75%' }
    @{ Row = 8;  B = 'MXToken';                           C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx';                           D = '1.060';      E = '14.07%' }
    @{ Row = 9;  B = 'LiechtensteinCryptoassetsExchange'; C = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx';    D = '0.1281';     E = '7.21%' }
    @{ Row = 10; B = 'WazirX';                            C = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx';                               D = '0.1877';     E = '2.28%' }
    @{ Row = 11; B = 'MandalaExchangeToken';              C = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx';            D = '0.09232';    E = '4.14%' }
    @{ Row = 12; B = 'BitrueCoin';                        C = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr';                          D = '0.04155';    E = '2.96%' }
    @{ Row = 13; B = 'BitMartToken';                      C = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx';                    D = '0.1039';     E = '-1.38%' }
    @{ Row = 14; B = 'BitForexToken';                     C = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf';                    D = '0.001305';   E = '1.80%' }
    @{ Row = 15; B = 'TigerCash';                         C = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch';                            D = '0.005715';   E = '-1.43%' }
    @{ Row = 16; B = 'UpBots';                            C = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt';                             D = '0.007409';   E = '1,889.69%' }
    @{ Row = 17; B = 'LEO';                               C = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo';                              D = '3.340';      E = '-0.07%' }
    @{ Row = 18; D = '2.373';    E = '-2.08%' }
    @{ Row = 19; D = '0.3430';   E = '3.61%' }
    @{ Row = 20; D = '8.023';    E = '0.38%' }
    @{ Row = 21; D = '0.1373';   E = '-3.21%' }
    @{ Row = 22; D = '0.2791';   E = '-6.89%' }
    @{ Row = 23; D = '0.04169';  E = '2.84%' }
    @{ Row = 24; D = '0.001270'; E = '0.33%' }
    @{ Row = 25; D = '0.004515'; E = '5.55%' }
    @{ Row = 26; D = '0.0001335';E = '8.59%' }
    @{ Row = 38; D = '0.02652';  E = '9.88%' }
    @{ Row = 39; D = '0.05373';  E = '3.06%' }
    @{ Row = 40; D = '0.005554'; E = '-14.51%' }
    @{ Row = 41; D = '0.007810'; E = '0.00%' }
    @{ Row = 42; D = '0.1389';   E = '4.25%' }
    @{ Row = 43; D = '0.007271'; E = '-2.38%' }
    @{ Row = 44; D = '0.008264'; E = '5.62%' }
    @{ Row = 45; D = '0.3028';   E = '-5.98%' }
    @{ Row = 46; D = '0.00006654'; E = '-0.63%' }
    @{ Row = 47; E = '-1.06%' }
    @{ Row = 48; D = '0.04780';  E = '3.30%' }
    @{ Row = 49; D = '0.003953'; E = '-5.84%' }
    @{ Row = 50; E = '-1.06%' }
    @{ Row = 51; E = '-1.06%' }
)

foreach ($update in $updates) {
    $row = $update.Row
    foreach ($col in 'B', 'C', 'D', 'E') {
        if ($update.ContainsKey($col)) {
            $cell = $ws.Range("$col$row")
            # Keep these as literal text (e.g. "305.91", "2.79%") instead of
            # letting Excel auto-convert them to numbers/percentages, since
            # the sheet stores every data cell as text.
            $cell.NumberFormat = "@"
            $cell.Value = $update[$col]
        }
    }
}
